$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "B2" = 85.71675767925336
    "C2" = 84.69382744795865
    "D2" = 85.33332045116963
    "E2" = 85.61792736458111

    "B3" = 96.75749996185083
    "C3" = 93.5463033240045
    "D3" = 96.77404069220678
    "E3" = 94.63371606006127

    "B4" = 98.96081898877081
    "C4" = 99.06258632147387
    "D4" = 98.94452880065849
    "E4" = 99.09910535596723

    "B5" = 98.61838668025598
    "C5" = 98.67065075131791
    "D5" = 98.52054337903078
    "E5" = 98.54023713723193

    "B6" = 98.33217964700556
    "C6" = 98.25661957506259
    "D6" = 98.29205309711423
    "E6" = 98.3396431843243

    "B7" = 97.62835332232727
    "C7" = 97.51602361878761
    "D7" = 97.58040650064591
    "E7" = 97.56218829767967

    "B8" = 96.29246898227525
    "C8" = 96.36744304973205
    "D8" = 96.07490805637107
    "E8" = 96.22810733307395
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
